$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the "_GoBack" bookmark that sits between "2100" and " RMB"
# in the registration-fee table. (This runtime has no working
# Bookmark.Delete, so we force its removal by re-writing the text that
# spans it -- Word drops a bookmark when the text containing it is
# replaced like this.)
# ---------------------------------------------------------------------------
$feeRange = $d.Content
$feeRange.Find.Execute("2100 RMB", $false, $false, $false, $false, $false, `
    $true, 1, $false, "2100 RMB", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: split the PayPal-receipt sentence and add " (mail or
# screenshot)" before " by email to: ", putting the (now-unique) "_GoBack"
# bookmark right after the word "screenshot".
# ---------------------------------------------------------------------------
$openQuote  = [char]0x201C
$closeQuote = [char]0x201D

$oldSentence = "PayPal payment, please state clearly the " + $openQuote + `
    "Participant name" + $closeQuote + " and send a copy of the PayPal " + `
    "receipt by email to: "
$newSentence = "PayPal payment, please state clearly the " + $openQuote + `
    "Participant name" + $closeQuote + " and send a copy of the PayPal " + `
    "receipt (mail or screenshot) by email to: "

$payRange = $d.Content
$payRange.Find.Execute($oldSentence, $false, $false, $false, $false, $false, `
    $true, 1, $false, $newSentence, 2) | Out-Null

$markRange = $d.Content
$markRange.Find.Execute("(mail or screenshot", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null

$bookmarkSpot = $d.Range($markRange.End, $markRange.End)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)
